# Apply cryptos list update (GitHub Actions scheduled refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-CellText {
    param($cellRange, [string]$text)
    if ($text -match '^[+-]?[0-9]+(\.[0-9]+)?$') {
        $cellRange.Value = "'" + $text
    } else {
        $cellRange.Value = $text
    }
}

Set-CellText $ws.Range("D2") "69.448.82"
Set-CellText $ws.Range("E2") "  -1.00%  "
Set-CellText $ws.Range("D3") "3.544.17"
Set-CellText $ws.Range("E3") "  -1.72%  "
Set-CellText $ws.Range("E4") "  +0.12%  "
Set-CellText $ws.Range("D5") "197.45"
Set-CellText $ws.Range("E5") "  +0.75%  "
Set-CellText $ws.Range("D6") "585.17"
Set-CellText $ws.Range("E6") "  -3.03%  "
Set-CellText $ws.Range("E7") "  -2.35%  "
Set-CellText $ws.Range("D8") "0.999"
Set-CellText $ws.Range("E9") "  -1.86%  "
Set-CellText $ws.Range("D10") "0.630"
Set-CellText $ws.Range("E10") "  -2.74%  "
Set-CellText $ws.Range("E11") "  -3.41%  "
Set-CellText $ws.Range("E12") "  -5.91%  "
Set-CellText $ws.Range("D13") "9.27"
Set-CellText $ws.Range("E13") "  -3.23%  "
Set-CellText $ws.Range("D14") "4.109.14"
Set-CellText $ws.Range("E14") "  -1.73%  "
Set-CellText $ws.Range("D15") "665.13"
Set-CellText $ws.Range("E15") "  +12.36%  "
Set-CellText $ws.Range("D16") "69.595.38"
Set-CellText $ws.Range("E16") "  -0.99%  "
Set-CellText $ws.Range("D17") "3.555.93"
Set-CellText $ws.Range("E17") "  -1.24%  "
Set-CellText $ws.Range("D18") "12.45"
Set-CellText $ws.Range("E18") "  -5.01%  "
Set-CellText $ws.Range("D19") "18.55"
Set-CellText $ws.Range("E19") "  -3.22%  "
Set-CellText $ws.Range("E20") "  -0.77%  "
Set-CellText $ws.Range("E21") "  -2.66%  "
Set-CellText $ws.Range("D22") "18.29"
Set-CellText $ws.Range("E22") "  +3.33%  "
Set-CellText $ws.Range("D23") "5.29"
Set-CellText $ws.Range("E23") "  +2.26%  "
Set-CellText $ws.Range("D24") "105.51"
Set-CellText $ws.Range("E24") "  +3.48%  "
Set-CellText $ws.Range("E25") "  -4.93%  "
Set-CellText $ws.Range("E26") "  -3.73%  "
Set-CellText $ws.Range("D27") "10.19"
Set-CellText $ws.Range("E27") "  -5.70%  "
Set-CellText $ws.Range("E28") "  +0.74%  "
Set-CellText $ws.Range("D29") "33.52"
Set-CellText $ws.Range("E29") "  -1.37%  "
Set-CellText $ws.Range("E30") "  -7.28%  "
Set-CellText $ws.Range("D31") "6.85"
Set-CellText $ws.Range("E31") "  -3.91%  "
Set-CellText $ws.Range("D32") "11.86"
Set-CellText $ws.Range("E32") "  -3.60%  "
Set-CellText $ws.Range("E33") "  -5.03%  "
Set-CellText $ws.Range("D34") "61.90"
Set-CellText $ws.Range("E34") "  -2.21%  "
Set-CellText $ws.Range("D35") "3.785.21"
Set-CellText $ws.Range("E35") "  -3.91%  "
Set-CellText $ws.Range("B36") "PEPE"
Set-CellText $ws.Range("C36") "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
Set-CellText $ws.Range("D36") "0.0₃0815"
Set-CellText $ws.Range("E36") "  -8.96%  "
Set-CellText $ws.Range("B37") "Stacks"
Set-CellText $ws.Range("C37") "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
Set-CellText $ws.Range("D37") "3.73"
Set-CellText $ws.Range("E37") "  +5.34%  "
Set-CellText $ws.Range("D38") "1.00"
Set-CellText $ws.Range("E38") "  +0.00%  "
Set-CellText $ws.Range("D39") "503.33"
Set-CellText $ws.Range("E39") "  -4.69%  "
Set-CellText $ws.Range("E40") "  -6.38%  "
Set-CellText $ws.Range("D41") "0.373"
Set-CellText $ws.Range("E41") "  -4.65%  "
Set-CellText $ws.Range("D42") "0.134"
Set-CellText $ws.Range("E42") "  +0.44%  "
Set-CellText $ws.Range("D43") "34.70"
Set-CellText $ws.Range("E43") "  -6.41%  "
Set-CellText $ws.Range("D44") "0.0453"
Set-CellText $ws.Range("E44") "  -0.47%  "
Set-CellText $ws.Range("D45") "2.89"
Set-CellText $ws.Range("E45") "  +1.13%  "
Set-CellText $ws.Range("E46") "  -1.13%  "
Set-CellText $ws.Range("D47") "0.136"
Set-CellText $ws.Range("E47") "  -3.23%  "
Set-CellText $ws.Range("E49") "  -3.41%  "
Set-CellText $ws.Range("D50") "1.80"
Set-CellText $ws.Range("E50") "  +21.56%  "
Set-CellText $ws.Range("D51") "2.70"
Set-CellText $ws.Range("E51") "  +62.29%  "
